$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.971.72"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").Value = "1.701.41"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3986"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4021"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.459"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.008"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08783"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.433"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001348"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.919"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("D17").Value = "1.705.67"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07228"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.222"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").Value = "24.972.72"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.909"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "150.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.252"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +20.86%  "

$ws.Range("D33").Value = "1.893.40"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08539"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.035"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.130"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2884"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09729"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8181"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.472"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.650"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7330"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.250"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.43%  "

$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.410"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.41%  "

